$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10
$ws.Range("B8").Value = 27
$ws.Range("B10").Value = 99
$ws.Range("B11").Value = 10
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 31
$ws.Range("B15").Value = 4
$ws.Range("B23").Value = 2
$ws.Range("B24").Value = 0
$ws.Range("B28").Value = 0
$ws.Range("B29").Value = 56
$ws.Range("B30").Value = 1
$ws.Range("B31").Value = 242
$ws.Range("B32").Value = 22
$ws.Range("B33").Value = 17
$ws.Range("B35").Value = 2
